$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.082.53'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '3.090.94'
$ws.Range('E3').Value = '  -2.44%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '234.57'
$ws.Range('E5').Value = '  +8.48%  '
$ws.Range('D6').Value = '620.15'
$ws.Range('E6').Value = '  -1.25%  '
$ws.Range('E7').Value = '  -13.04%  '
$ws.Range('D8').Value = '0.359'
$ws.Range('E8').Value = '  -3.77%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '3.088.96'
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('D11').Value = '0.717'
$ws.Range('E11').Value = '  -6.76%  '
$ws.Range('D12').Value = '0.197'
$ws.Range('E12').Value = '  -3.12%  '
$ws.Range('D13').Value = '0.0000247'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').Value = '35.09'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = '89.875.80'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('D16').Value = '5.38'
$ws.Range('E16').Value = '  -6.88%  '
$ws.Range('E17').Value = '  -2.76%  '
$ws.Range('D18').Value = '3.085.04'
$ws.Range('E18').Value = '  -3.71%  '
$ws.Range('D19').Value = '3.82'
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('D20').Value = '0.0000211'
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('D21').Value = '13.82'
$ws.Range('E21').Value = '  -6.39%  '
$ws.Range('D22').Value = '432.09'
$ws.Range('E22').Value = '  -10.02%  '
$ws.Range('D23').Value = '5.46'
$ws.Range('E23').Value = '  +2.82%  '
$ws.Range('D24').Value = '8.79'
$ws.Range('E24').Value = '  -4.21%  '
$ws.Range('D25').Value = '5.58'
$ws.Range('E25').Value = '  -3.34%  '
$ws.Range('D26').Value = '86.10'
$ws.Range('E26').Value = '  -11.05%  '
$ws.Range('D27').Value = '11.80'
$ws.Range('E27').Value = '  -4.99%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').Value = '9.09'
$ws.Range('E30').Value = '  -2.43%  '
$ws.Range('E31').Value = '  +3.12%  '
$ws.Range('E32').Value = '  -4.83%  '
$ws.Range('E33').Value = '  -3.95%  '
$ws.Range('D34').Value = '25.62'
$ws.Range('E34').Value = '  -9.52%  '
$ws.Range('E35').Value = '  +3.50%  '
$ws.Range('D36').Value = '3.72'
$ws.Range('E36').Value = '  +2.23%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '7.11'
$ws.Range('E37').Value = '  +1.83%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '498.12'
$ws.Range('E38').Value = '  -5.35%  '
$ws.Range('E39').Value = '  -3.00%  '
$ws.Range('E40').Value = '  -3.51%  '
$ws.Range('D41').Value = '3.65'
$ws.Range('E41').Value = '  +57.46%  '
$ws.Range('D42').Value = '0.0868'
$ws.Range('E42').Value = '  -4.26%  '
$ws.Range('D43').Value = '22.10'
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = '0.399'
$ws.Range('E45').Value = '  -5.50%  '
$ws.Range('E46').Value = '  -6.32%  '
$ws.Range('D47').Value = '0.680'
$ws.Range('E47').Value = '  -4.03%  '
$ws.Range('D48').Value = '150.80'
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('D49').Value = '44.39'
$ws.Range('E49').Value = '  -2.26%  '
$ws.Range('E50').Value = '  -4.85%  '
$ws.Range('D51').Value = '0.999'
$ws.Range('E51').Value = '  -0.13%  '
